$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SKU / url values for the two data rows
$ws.Range("A2").Value = "118032724"
$ws.Range("B2").Value = "https://imagedelivery.net/4fYuQyy-r8_rpBpcY7lH_A/falabellaPE/118032724_01/w=800,h=800,fit=pad"

$ws.Range("A3").Value = "127356818"
$ws.Range("B3").Value = "https://imagedelivery.net/4fYuQyy-r8_rpBpcY7lH_A/falabellaPE/127356818_01/w=800,h=800,fit=pad|https://imagedelivery.net/4fYuQyy-r8_rpBpcY7lH_A/falabellaPE/127356818_02/w=800,h=800,fit=pad|https://imagedelivery.net/4fYuQyy-r8_rpBpcY7lH_A/falabellaPE/127356818_03/w=800,h=800,fit=pad|https://imagedelivery.net/4fYuQyy-r8_rpBpcY7lH_A/falabellaPE/127356818_04/w=800,h=800,fit=pad"

# Restore column A to its default (non-custom) width
$ws.Columns.Item(1).ColumnWidth = 10.81

# Move the active selection to D10
$ws.Range("D10").Select()
